$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns in front of the existing "B:D" columns.
# This pushes the current B,C,D,E (Jun_17, Jun_15, Jun_13, Jun_10) out to E,F,G,H
# and leaves three brand-new, empty columns at B, C, D.
$ws.Columns("B:D").Insert()

# Give every data column (C through H) the same explicit 8-character width
# that the sheet already used for the old C/D/E columns.
$ws.Columns("C:H").ColumnWidth = 7.166666666666667

# New week headers in row 1 for the freshly inserted columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# The insert left B2:D27 blank for every existing analyst row - these weeks
# have no rating action recorded, so mark them "UN" like the rest of the sheet.
$ws.Range("B2:D27").Value = "UN"

# Two new analyst/firm rows appended at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
